# Apply the "Add data for 2022-12-08" update:
# - Rename sheet (and its header label) from "...11-29" to "...11-30"
# - Update December (row 12) and Total (row 14) figures in column I

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the worksheet tab
$ws.Name = "Through 2022-11-30"

# Update the column header label (shared string used by I1)
$ws.Range("I1").Value = "2022 (through 11-30)"

# Update the December and Total figures in column I
$ws.Range("I12").Value = 118
$ws.Range("I14").Value = 1516
